$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update D2 email address
$ws.Range("D2").Value = "riya.ingale14@gmail.com"

# Update O2 and P2 dates
$ws.Range("O2").Value = "20 /01 /21"
$ws.Range("P2").Value = "22/ 02/ 21"

# Update T2 value
$ws.Range("T2").Value = "Yes"

# Add new header for column U
$ws.Range("U1").Value = "Completition Certificate"

# Add hyperlink formula in U2
$ws.Range("U2").Formula = '=HYPERLINK("http://127.0.0.1:5000/downloadcompletioncert/1", "Download Cert")'

# Delete rows 3 and 4 (Samiksha Pansare and Deep)
$ws.Rows("3:4").Delete()
